$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# 1. Add a new column to the table. The engine always appends it physically at
#    the end (column V) and gives it a fresh id - this mirrors Excel's own
#    behaviour of assigning inserted table columns a brand-new id.
$newCol = $tbl.ListColumns.Add()

# 2. Re-point the header row so "Longitudinal" ends up physically in column D,
#    pushing the old D.. headers one slot to the right ("Treatment variable"
#    -> E, "Column3" renamed to "Questions" -> F, Column4..Column19 -> G..V).
$ws.Range("D1").Value = "Longitudinal"
$ws.Range("E1").Value = "Treatment variable"
$ws.Range("F1").Value = "Questions"
$ws.Range("G1").Value = "Column4"
$ws.Range("H1").Value = "Column5"
$ws.Range("I1").Value = "Column6"
$ws.Range("J1").Value = "Column7"
$ws.Range("K1").Value = "Column8"
$ws.Range("L1").Value = "Column9"
$ws.Range("M1").Value = "Column10"
$ws.Range("N1").Value = "Column11"
$ws.Range("O1").Value = "Column12"
$ws.Range("P1").Value = "Column13"
$ws.Range("Q1").Value = "Column14"
$ws.Range("R1").Value = "Column15"
$ws.Range("S1").Value = "Column16"
$ws.Range("T1").Value = "Column17"
$ws.Range("U1").Value = "Column18"
$ws.Range("V1").Value = "Column19"

# 3. PMID for the BIDMC-FMT row.
$ws.Range("B2").Value = 27542133

# 4. Shift the one pre-existing "Treatment variable" value (row 2's
#    "baseline") one column right into the new column E, then populate the
#    new Longitudinal column (D) and the new Questions column (F) for row 2.
$ws.Range("E2").Value = "baseline"
$ws.Range("D2").Value = "Yes according to paper. Samples before and after FMT"
$ws.Range("F2").Value = "are samples baseline or???"

# 5. Row 3 gets a new "Treatment variable" note.
$ws.Range("E3").Value = "No treatment"

# 6. Row 15's old "Treatment variable" value ("   ") shifts right into E15;
#    clear it back out of D15 (formatting for D15 is reapplied below).
$ws.Range("E15").Value = "   "
$ws.Range("D15").ClearContents()

# 7. Give every Longitudinal cell (D2:D16) the same formatting as the
#    adjacent "n" cell on its row, matching the per-row font used in column C
#    (rows 12-13 use the Helvetica Neue style, the rest use Arial).
for ($row = 2; $row -le 16; $row++) {
    $src = $ws.Cells.Item($row, 3)
    $dst = $ws.Cells.Item($row, 4)
    $src.Copy()
    $dst.PasteSpecial(-4122)
}
$ws.Range("A1").Select()
$excel.CutCopyMode = 0

# 8. Match the saved selection from the authored workbook.
$ws.Range("D3").Select()
